$d = $word.ActiveDocument
$wns = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# 1) Paragraph "The purpose of the document: ... 1.0." -> drop the trailing
#    _GoBack bookmark (it moves to the next paragraph, see step 2).
$target1 = $d.Content.Find
$target1.Execute("specification of a centrifuge protocol") | Out-Null
$p1 = $target1.Parent.Paragraphs(1).Range
$xml1 = '<w:p ' + $wns + ' w:rsidR="00783D22" w:rsidRDefault="00783D22" w:rsidP="00783D22">' +
  '<w:r><w:t>Th</w:t></w:r>' +
  '<w:r w:rsidR="003A69B0"><w:t>e</w:t></w:r>' +
  '<w:r><w:t xml:space="preserve"> </w:t></w:r>' +
  '<w:r w:rsidR="003A69B0"><w:t xml:space="preserve">purpose of the document: </w:t></w:r>' +
  '<w:r><w:t xml:space="preserve">specification of a centrifuge protocol </w:t></w:r>' +
  '<w:r w:rsidR="00200A94"><w:t>ordinal</w:t></w:r>' +
  '<w:r><w:t xml:space="preserve"> 1.0.</w:t></w:r>' +
  '</w:p>'
$p1.InsertXML($xml1) | Out-Null

# 2) Paragraph "The protocol ordinal 1.0 is applied to a single dataset...."
#    rewritten with the new sentences about "Values in each column" and the
#    gone proofErr markers; the _GoBack bookmark now lives at its end.
$target2 = $d.Content.Find
$target2.Execute("The protocol ") | Out-Null
$p2 = $target2.Parent.Paragraphs(1).Range
$xml2 = '<w:p ' + $wns + ' w:rsidR="00783D22" w:rsidRDefault="00783D22" w:rsidP="00783D22">' +
  '<w:r><w:t xml:space="preserve">The protocol </w:t></w:r>' +
  '<w:r w:rsidR="00200A94"><w:t xml:space="preserve">ordinal </w:t></w:r>' +
  '<w:r><w:t>1.0 is applied to a single dataset. If a data owner wish</w:t></w:r>' +
  '<w:r><w:t>es</w:t></w:r>' +
  '<w:r><w:t xml:space="preserve"> to apply it to a database containing more than one dataset, the database must be converted to a single dataset before applying the protocol. </w:t></w:r>' +
  '<w:r><w:t>Values in e</w:t></w:r>' +
  '<w:r><w:t xml:space="preserve">ach column </w:t></w:r>' +
  '<w:r><w:t xml:space="preserve">are replaced with integer ordinal values. If there are N unique values in the column, then the minimal value is replaced with a 1, the second minimal is replaced with a 2, etc. The maximal value is replaced with N. </w:t></w:r>' +
  '<w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/>' +
  '</w:p>'
$p2.InsertXML($xml2) | Out-Null

# 3) Table cell "SUBJID" -> drop the stray lastRenderedPageBreak marker.
$target3 = $d.Content.Find
$target3.Execute("SUBJID") | Out-Null
$p3 = $target3.Parent.Paragraphs(1).Range
$xml3 = '<w:p ' + $wns + ' w:rsidR="00856A9D" w:rsidRDefault="00856A9D" w:rsidP="00207842">' +
  '<w:r><w:t>SUBJID</w:t></w:r>' +
  '</w:p>'
$p3.InsertXML($xml3) | Out-Null
